$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Online kein Bestand" note for the Betty Bossi Pfaffenhut row (row 305)
$ws.Cells.Item(305, 13).Value = "Betty Bossi Pfaffenhut 2x  100g - Online kein Bestand 4.60 Schweizer Franken"

# Refresh the scrape timestamp (column O) for every data row (2..388) to the new crawl time
for ($r = 2; $r -le 388; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-02-22 20:49:40"
}

# The product "Naturaplan Bio Appenzeller Herzbiber 3x76g" (row 369) is no longer present
# in this crawl, so remove the whole row - everything below shifts up by one.
$ws.Rows(369).Delete()
